$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: new bitácora entry, same formatting layout as row 2
$ws.Range("A3").Value = "Se realizó la interfaz de gestión de clientes, se creó el repositorio en github y se creó el primer commit."

$ws.Range("B2").Copy()
$ws.Range("B3:C3").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("B3").Value = 42863.416666666664
$ws.Range("C3").Value = 42862.541666666664
$ws.Range("D3").Value = 0.125

$ws.Rows("3").RowHeight = 30

# Move the active selection to D4
$ws.Range("D4").Select()
